# Applies odds-update changes to the "Jogos da Semana" worksheet
# (rows 2-7 of Sheet1), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("I2").Value = 5.25
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("AC2").Value = 6.5
$ws.Range("AH2").Value = 23

# Row 3
$ws.Range("G3").Value = 2.7
$ws.Range("I3").Value = 2.9
$ws.Range("J3").Value = 3.6
$ws.Range("W3").Value = 6.5
$ws.Range("Z3").Value = 29
$ws.Range("AC3").Value = 6
$ws.Range("AG3").Value = 6.5
$ws.Range("AH3").Value = 12
$ws.Range("AJ3").Value = 29
$ws.Range("AN3").Value = 4.5
$ws.Range("AS3").Value = 351
$ws.Range("BD3").Value = 126

# Row 4
$ws.Range("G4").Value = 3.6
$ws.Range("I4").Value = 2.25
$ws.Range("J4").Value = 4.33
$ws.Range("L4").Value = 3.1
$ws.Range("W4").Value = 8.5
$ws.Range("X4").Value = 17
$ws.Range("Z4").Value = 41
$ws.Range("AC4").Value = 6.5
$ws.Range("AH4").Value = 9.5
$ws.Range("AJ4").Value = 21
$ws.Range("AK4").Value = 21
$ws.Range("AQ4").Value = 81
$ws.Range("AR4").Value = 126
$ws.Range("AW4").Value = 4
$ws.Range("AX4").Value = 13
$ws.Range("BB4").Value = 251

# Row 5
$ws.Range("G5").Value = 1.57
$ws.Range("H5").Value = 3.9
$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 2.2
$ws.Range("L5").Value = 6.5
$ws.Range("U5").Value = 2.2
$ws.Range("V5").Value = 1.62
$ws.Range("W5").Value = 5.5
$ws.Range("X5").Value = 6.5
$ws.Range("Z5").Value = 11
$ws.Range("AC5").Value = 9
$ws.Range("AD5").Value = 7.5
$ws.Range("AF5").Value = 81
$ws.Range("AH5").Value = 29
$ws.Range("AK5").Value = 51
$ws.Range("AO5").Value = 8
$ws.Range("AQ5").Value = 26
$ws.Range("AW5").Value = 7.5

# Row 6
$ws.Range("G6").Value = 1.9
$ws.Range("I6").Value = 4.5
$ws.Range("J6").Value = 2.63
$ws.Range("L6").Value = 5
$ws.Range("M6").Value = 1.08
$ws.Range("N6").Value = 8
$ws.Range("Z6").Value = 15
$ws.Range("AA6").Value = 17
$ws.Range("AD6").Value = 6.5
$ws.Range("AG6").Value = 10
$ws.Range("AJ6").Value = 51
$ws.Range("AM6").Value = 1000
$ws.Range("AX6").Value = 26
$ws.Range("AY6").Value = 41
$ws.Range("AZ6").Value = 101

# Row 7
$ws.Range("G7").Value = 1.9
$ws.Range("I7").Value = 4.75
$ws.Range("J7").Value = 2.63
$ws.Range("K7").Value = 1.95
$ws.Range("L7").Value = 5
$ws.Range("M7").Value = 1.11
$ws.Range("N7").Value = 6.5
$ws.Range("AD7").Value = 6
$ws.Range("AH7").Value = 21
$ws.Range("AN7").Value = 3.75
$ws.Range("AO7").Value = 11
$ws.Range("AX7").Value = 26
